$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.946.08"
$ws.Range("E2").Value = "  +7.97%  "
$ws.Range("D3").Value = "1.826.71"
$ws.Range("E3").Value = "  +5.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.19"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4937"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.81"
$ws.Range("E8").Value = "  +6.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2786"
$ws.Range("E9").Value = "  +7.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06409"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("D11").Value = "1.808.88"
$ws.Range("E11").Value = "  +4.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.76"
$ws.Range("E12").Value = "  +5.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07058"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6456"
$ws.Range("E14").Value = "  +6.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.30"
$ws.Range("E15").Value = "  +9.54%  "
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "28.978.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9983"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007321"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.26"
$ws.Range("E21").Value = "  +7.77%  "
$ws.Range("D22").Value = "2.040.73"
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.572"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.859"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.48"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "129.38"
$ws.Range("E27").Value = "  +21.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.44"
$ws.Range("E28").Value = "  +7.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.890"
$ws.Range("E29").Value = "  +5.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.413"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.134"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08369"
$ws.Range("E32").Value = "  +5.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.794"
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04932"
$ws.Range("E34").Value = "  +8.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.102"
$ws.Range("E35").Value = "  +10.10%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6733"
$ws.Range("E36").Value = "  +8.69%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.296"
$ws.Range("E38").Value = "  +14.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.702"
$ws.Range("E39").Value = "  +10.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9479"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.171"
$ws.Range("E41").Value = "  +9.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.54"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4084"
$ws.Range("E45").Value = "  +6.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.225"
$ws.Range("E46").Value = "  +6.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1225"
$ws.Range("E47").Value = "  +5.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05525"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.211"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.68"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.306"
$ws.Range("E51").Value = "  +4.83%  "
